$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update pressure readings (column B) that changed in the refreshed upload
$ws.Range("B6").Value = 128
$ws.Range("B7").Value = 133
$ws.Range("B8").Value = 137
$ws.Range("B36").Value = 134
$ws.Range("B39").Value = 139
$ws.Range("B49").Value = 133

# Reflect where the author had scrolled/selected when the file was saved
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("B49").Select()
